$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row -> Value column gets "false".
# A plain string assignment of "false" gets auto-coerced to a native boolean
# by the engine, so instead write it as a formula string result and then
# paste-special as values to collapse it back into a genuine text cell
# (matching the shared-string "false" in the target workbook) while keeping
# the existing cell style intact.
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# "Date" row -> refresh the generation timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# "Description" row -> Value column gets the CodeSystem description text.
$ws.Range("B17").Value = "Standards and norms used for CRF classification"
